# Apply corrections to part 3/4 motif results sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C (Triplet) values for rows 2-20 (C2:C20), replacing the old
# motif triplet list with the corrected one.
$triplets = @(
    "[57, 47, 43]",
    "[57, 47, 46]",
    "[61, 57, 46]",
    "[62, 57, 46]",
    "[50, 57, 43]",
    "[57, 43, 64]",
    "[62, 57, 64]",
    "[48, 57, 43]",
    "[49, 57, 43]",
    "[56, 57, 43]",
    "[57, 58, 16]",
    "[57, 58, 46]",
    "[61, 56, 57]",
    "[61, 57, 58]",
    "[62, 57, 58]",
    "[49, 57, 56]",
    "[50, 49, 57]",
    "[48, 49, 57]",
    "[48, 56, 57]"
)

# Column B (counts) values for rows 2-20.
$counts = @(38, 38, 38, 38, 38, 38, 38, 46, 46, 46, 46, 46, 108, 108, 108, 110, 110, 238, 238)

# Rows 12-20 are brand new; copy the formatting already used by column A
# in the existing data rows (e.g. A11) before filling in their values.
$ws.Range("A11").Copy() | Out-Null
$ws.Range("A12:A20").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

for ($i = 0; $i -lt $triplets.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $i
    $ws.Cells.Item($row, 2).Value = $counts[$i]
    $ws.Cells.Item($row, 3).Value = $triplets[$i]
}
